$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite B1:E3 with the values from N1:N3, O1:O3, AM1:AM3, AN1:AN3
# (subjects 15 and 16, CON and STR groups) replacing the old placeholder
# "Subj" header columns (B-E) and leg/ROM data.

$ws.Range("B1").Value = $ws.Range("N1").Value()
$ws.Range("C1").Value = $ws.Range("O1").Value()
$ws.Range("D1").Value = $ws.Range("AM1").Value()
$ws.Range("E1").Value = $ws.Range("AN1").Value()

$ws.Range("B2").Value = $ws.Range("N2").Value()
$ws.Range("C2").Value = $ws.Range("O2").Value()
$ws.Range("D2").Value = $ws.Range("AM2").Value()
$ws.Range("E2").Value = $ws.Range("AN2").Value()

$ws.Range("B3").Value = $ws.Range("N3").Value()
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = $ws.Range("AM3").Value()
$ws.Range("E3").Value = $ws.Range("AN3").Value()

# Update the active selection to reflect the newly updated block
$ws.Range("B1:E3").Select()
